$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Cell values.
#    Text cells are written in a specific order so the generated shared
#    strings table (xl/sharedStrings.xml) lines up with the original file:
#    Peso(g)=0, Nome=1, Lateral (Al)=2, Cobertura (A definir)=3,
#    Chassi (PETG)=4, motor arma=5, bateria=6, receptor=7, ESC=8,
#    Suporte motor mov(Al)=9, Arma (Aço + Al)=10, Rolamentos=11,
#    Correia + Polias=12, Placa de controle=13, motores mov=14,
#    Fios + conectores=15, Parafusos=16, TOTAL=17, Rodas=18,
#    Mancal  arma (Pol)=19
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Peso(g)"
$ws.Range("A1").Value = "Nome"
$ws.Range("A2").Value = "Lateral (Al)"
$ws.Range("A3").Value = "Cobertura (A definir)"
$ws.Range("A4").Value = "Chassi (PETG)"
$ws.Range("A12").Value = "motor arma"
$ws.Range("A13").Value = "bateria"
$ws.Range("A14").Value = "receptor"
$ws.Range("A15").Value = "ESC"
$ws.Range("A5").Value = "Suporte motor mov(Al)"
$ws.Range("A7").Value = "Arma (Aço + Al)"
$ws.Range("A8").Value = "Rolamentos"
$ws.Range("A9").Value = "Correia + Polias"
$ws.Range("A16").Value = "Placa de controle"
$ws.Range("A10").Value = "motores mov"
$ws.Range("A17").Value = "Fios + conectores"
$ws.Range("A18").Value = "Parafusos"
$ws.Range("A19").Value = "TOTAL"
$ws.Range("A11").Value = "Rodas"
$ws.Range("A6").Value = "Mancal  arma (Pol)"

# Numeric weights (grams)
$ws.Range("B2").Value = 140
$ws.Range("B3").Value = 150
$ws.Range("B4").Value = 180
$ws.Range("B5").Value = 50
$ws.Range("B6").Value = 50
$ws.Range("B7").Value = 223
$ws.Range("B8").Value = 30
$ws.Range("B9").Value = 100
$ws.Range("B10").Value = 170
$ws.Range("B11").Value = 68
$ws.Range("B12").Value = 50
$ws.Range("B13").Value = 86
$ws.Range("B14").Value = 8
$ws.Range("B15").Value = 25
$ws.Range("B16").Value = 50
$ws.Range("B17").Value = 50
$ws.Range("B18").Value = 120

# Total (sum) row
$ws.Range("B19").Formula = "=SUM(B2:B18)"

# ---------------------------------------------------------------------------
# 2) Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 24
$ws.Columns.Item(2).ColumnWidth = 10.6

# ---------------------------------------------------------------------------
# 3) Cell formatting. Applied in an order that reproduces the original
#    cellXfs table: (1) header cell border+center, (2) plain border used
#    throughout the table, (3) total-row fill+border, (4) the "Rodas" row
#    left/right-only border, (5) red font for highlighted rows.
# ---------------------------------------------------------------------------

# (1) Header name cell: box border + centered text
$ws.Range("A1").Borders.LineStyle = 1
$ws.Range("A1").HorizontalAlignment = -4108

# (2) Plain box-bordered cells used across most of the table
$ws.Range("B1").Borders.LineStyle = 1
$ws.Range("A4:B4").Borders.LineStyle = 1
$ws.Range("A6:B7").Borders.LineStyle = 1
$ws.Range("A10:B10").Borders.LineStyle = 1
$ws.Range("A12:B18").Borders.LineStyle = 1

# (3) Total row: yellow fill + box border
$ws.Range("A19:B19").Borders.LineStyle = 1
$ws.Range("A19:B19").Interior.Color = 65535

# (4) "Rodas" row: left & right border only (no top/bottom)
$ws.Range("A11").Borders.Item(7).LineStyle = 1
$ws.Range("A11").Borders.Item(10).LineStyle = 1
$ws.Range("A11").Borders.Item(8).LineStyle = -4142
$ws.Range("A11").Borders.Item(9).LineStyle = -4142

# (5) Highlighted rows: red font + box border
$ws.Range("A2:B3").Borders.LineStyle = 1
$ws.Range("A2:B3").Font.Color = 255
$ws.Range("A5:B5").Borders.LineStyle = 1
$ws.Range("A5:B5").Font.Color = 255
$ws.Range("A8:B9").Borders.LineStyle = 1
$ws.Range("A8:B9").Font.Color = 255

# ---------------------------------------------------------------------------
# 4) Selection shown when the workbook is (re)opened
# ---------------------------------------------------------------------------
$ws.Range("B17").Select()

Write-Output "done"
